$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 07:04"

# Tailandia (row 69) - updated covid numbers
$ws.Range("B69").Value = 3009
$ws.Range("C69").Value = 5
$ws.Range("D69").Value = 2794
$ws.Range("E69").Value = 159

# Guyana overtakes Monaco (rows 163/164 swap order with updated numbers)
$ws.Range("A163").Value = "Guyana"
$ws.Range("B163").Value = 97
$ws.Range("C163").Value = 3
$ws.Range("D163").Value = 35
$ws.Range("E163").Value = 52
$ws.Range("F163").Value = 5
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 10

$ws.Range("A164").Value = "Monaco"
$ws.Range("B164").Value = 96
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 82
$ws.Range("E164").Value = 10
$ws.Range("F164").Value = 1
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 4

# Belice and Nueva Caledonia swap order (rows 192/193)
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
